$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Widen columns I and K ---
$ws.Range("I1").ColumnWidth = 28.35
$ws.Range("K1").ColumnWidth = 40.5

# --- Insert two new rows at the top of the data (rows 3-4), pushing the
#     existing three data rows down to rows 5-7 ---
$ws.Rows("3:4").Insert()

# Clear the formatting the insert copied down from the header row, then
# restore the plain "integer" number format used by column A/D data cells.
$ws.Range("A3:K4").ClearFormats()
$ws.Range("A3").NumberFormat = "0"
$ws.Range("D3").NumberFormat = "0"
$ws.Range("A4").NumberFormat = "0"
$ws.Range("D4").NumberFormat = "0"

# --- Row 3: brand new walk-in entry (JOEL) ---
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "14-12-2025"
$ws.Range("C3").Value = "JOEL"
$ws.Range("D3").Value = 8086460359
$ws.Range("E3").Value = "20-01-2026"
$ws.Range("F3").Value = "Abdul Hadi Rafeeque"
$ws.Range("G3").Value = "Loss"
$ws.Range("H3").Value = "SIZE NOT SUITABLE"
$ws.Range("I3").Value = "SIZE TOO LARGE"
$ws.Range("J3").Value = "-"
$ws.Range("K3").Value = "SIZE ISSUE  ( 46 SIZE SUITE )"

# --- Row 4: brand new walk-in entry (zahid) ---
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "16-12-2025"
$ws.Range("C4").Value = "zahid"
$ws.Range("D4").Value = 7034372731
$ws.Range("E4").Value = "16-12-2025"
$ws.Range("F4").Value = "Abdul Hadi Rafeeque"
$ws.Range("G4").Value = "Loss"
$ws.Range("H4").Value = "PRODUCT"
$ws.Range("I4").Value = "PRODUCT NOT AVAILABLE"
$ws.Range("J4").Value = "-"
$ws.Range("K4").Value = "costomer will confirm tomorrow"

# --- Row 6 (previously row 4, "shamil"): fix the Staff value to match
#     the correction that was also applied when that row was first added ---
$ws.Range("F6").Value = "Abdul Hadi Rafeeque"
